# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") on the active sheet for rows 2-31 with the
# newly computed strikeout (K) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 3
    4  = 4
    5  = 3
    6  = 1
    7  = 2
    8  = 4
    9  = 6
    10 = 3
    11 = 1
    12 = 1
    13 = 6
    14 = 4
    15 = 4
    16 = 4
    17 = 4
    18 = 1
    19 = 0
    20 = 6
    21 = 4
    22 = 6
    23 = 5
    24 = 3
    25 = 5
    26 = 3
    27 = 6
    28 = 5
    29 = 5
    30 = 2
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
